# bug fix in Eduati data files
#
# Sheet1 ("SNUC2B_noCTRL_meas") only has real measurement data through
# row 44; rows 45-87 are leftover row-index stubs (column A only) that
# should not be part of the sheet. Trim them, and update the view state
# so Sheet1 (not Sheet3) is the active/selected sheet when the workbook
# is reopened.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the stray trailing rows (45-87) on Sheet1 - they only ever held
# a bare row counter in column A, no real data.
$ws1.Range("A45:A87").EntireRow.Delete() | Out-Null

# Sheet1 becomes the active sheet/tab (was Sheet3 before this fix; activating
# Sheet1 here also clears Sheet3's tabSelected flag).
$ws1.Activate() | Out-Null

# Restore the working selection on Sheet1.
$ws1.Range("F61").Select() | Out-Null
